$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.954.74"
$ws.Range("E2").Value = "  -2.05%  "

$ws.Range("D3").Value = "1.742.58"
$ws.Range("E3").Value = "  -0.31%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.0000"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.38%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "312.12"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -4.92%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.0000"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.37%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.5021"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +5.10%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3579"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +1.37%  "

$ws.Range("E9").Value = "  -0.50%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.07274"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -2.59%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "1.062"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -1.11%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.000"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -0.30%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "20.30"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -0.06%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.993"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -0.65%  "

$ws.Range("D15").Value = "1.745.94"
$ws.Range("E15").Value = "  +0.06%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "6.864"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -2.38%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "86.93"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -5.47%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.00001035"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -3.24%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.06404"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +0.20%  "

$ws.Range("E20").Value = "  -0.41%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "16.59"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -0.44%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.737"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -0.26%  "

$ws.Range("D23").Value = "27.036.22"
$ws.Range("E23").Value = "  -1.94%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "11.28"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +2.65%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.047"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -5.05%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "154.70"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -4.46%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "19.98"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +0.25%  "

$ws.Range("D28").Value = "1.943.13"
$ws.Range("E28").Value = "  -0.24%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.140"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -2.54%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "121.53"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -0.09%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.053"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +0.41%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.09530"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +2.09%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "3.576"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -1.34%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "5.391"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -1.51%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.02208"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -1.67%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.05896"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -0.64%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "11.11"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -1.89%  "

$ws.Range("E38").Value = "  +0.04%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.2000"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -2.27%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "4.763"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -1.38%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.6035"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -0.75%  "

$ws.Range("E42").Value = "  -0.41%  "

$ws.Range("E43").Value = "  -5.39%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "7.521"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -2.48%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "12.85"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -0.36%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "3.602"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -3.43%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.5656"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -0.84%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "119.92"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -2.14%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.862"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -2.24%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.108"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -2.27%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.06671"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -1.15%  "
